$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11 (tc_tourtpk_09 / content manager dashboard) ---
# C11 text is unchanged; add new "Steps to perform" text in D11.
$ws.Range("D11").WrapText = $true
$ws.Range("D11").Value = "1. go to the login page`n2. enter the valid data to log in`n3. able to create a new place`n4. able to view and delete place`n5. able to create a new plan`n6. able to view and delete plan`n7. able to view and delete preference`n8. logout."
$ws.Rows(11).RowHeight = 150

# --- Row 12 (tc_tourtpk_10 / admin dashboard) ---
# C12 scenario description updated to mention "partner" too.
$ws.Range("C12").Value = "admin  view, and deletes a details of places, recommended plan, hotels, and packages. Able to view and delete users and partner"
$ws.Range("D12").WrapText = $true
$ws.Range("D12").Value = "1. go to the login page`n2. enter the valid data to log in`n3. able to view and delete place`n4. able to view and delete packages`n5. able to view and delete plan`n6. able to view and delete hotels`n7. able to view and delete users `n8. logout."
$ws.Rows(12).RowHeight = 180

# --- Row 13 (tc_tourtpk_11 / hotel manager dashboard) ---
# C13 text is unchanged; add new "Steps to perform" text in D13.
$ws.Range("D13").WrapText = $true
$ws.Range("D13").Value = "1. go to the login page`n2. enter the valid data to log in`n3. able to create a new room`n4.able to update its profile`n5. able to view and delete hotels`n6. able to view and delete booking `n7. logout."
$ws.Rows(13).RowHeight = 135

# --- Row 14 (tc_tourtpk_12 / tour operator dashboard) ---
# C14 text is unchanged; add new "Steps to perform" text in D14.
$ws.Range("D14").WrapText = $true
$ws.Range("D14").Value = "1. go to the login page`n2. enter the valid data to log in`n3. able to create a new packages`n4. able to update its profile`n5. able to view and delete packages`n6. logout."
$ws.Rows(14).RowHeight = 120

# --- View state: scroll down so row 13 is at the top, select D17 ---
[void]$ws.Range("D17").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
